# Update lecture 1 materials for 2024
$p = $ppt.ActivePresentation

# 1. Update the date shown on the title slide.
$s1 = $p.Slides.Item(1)
$dateShape = $s1.Shapes.Item(3)
$dateParagraph = $dateShape.TextFrame.TextRange.Paragraphs(3)
$dateParagraph.Runs(1).Text = "January 4, 2024"

# 2. Swap the deck's theme palette ("Simple Light" -> "Default"
#    Highcharts-style palette) on the slide master's theme.
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme
$tcs.Item(1).RGB  = 0        # dk1      000000
$tcs.Item(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Item(3).RGB  = 5800213  # dk2      158158
$tcs.Item(4).RGB  = 15987699 # lt2      F3F3F3
$tcs.Item(5).RGB  = 13077765 # accent1  058DC7
$tcs.Item(6).RGB  = 3322960  # accent2  50B432
$tcs.Item(7).RGB  = 1791725  # accent3  ED561B
$tcs.Item(8).RGB  = 61421    # accent4  EDEF00
$tcs.Item(9).RGB  = 15059748 # accent5  24CBE5
$tcs.Item(10).RGB = 7529828  # accent6  64E572
$tcs.Item(11).RGB = 13369378 # hlink    2200CC
$tcs.Item(12).RGB = 9116245  # folHlink 551A8B
